$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'23.113.17"
$ws.Range('E2').Value = '  -3.24%  '
$ws.Range('D3').Value = "'1.603.00"
$ws.Range('E3').Value = '  -3.26%  '
$ws.Range('D4').Value = "'1.000"
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('E5').Value = '  -0.02%  '
$ws.Range('D6').Value = "'301.40"
$ws.Range('E6').Value = '  -3.19%  '
$ws.Range('D7').Value = "'0.3777"
$ws.Range('E7').Value = '  -3.03%  '
$ws.Range('D8').Value = "'0.3654"
$ws.Range('E8').Value = '  -4.37%  '
$ws.Range('D9').Value = "'49.85"
$ws.Range('E9').Value = '  -3.29%  '
$ws.Range('D10').Value = "'1.268"
$ws.Range('E10').Value = '  -6.33%  '
$ws.Range('D11').Value = "'0.08150"
$ws.Range('E11').Value = '  -4.00%  '
$ws.Range('D12').Value = "'1.000"
$ws.Range('E12').Value = '  -0.05%  '
$ws.Range('E13').Value = '  -4.50%  '
$ws.Range('D14').Value = "'6.592"
$ws.Range('D15').Value = "'0.00001260"
$ws.Range('E15').Value = '  -4.12%  '
$ws.Range('D16').Value = "'7.389"
$ws.Range('D17').Value = "'1.599.29"
$ws.Range('E17').Value = '  -3.40%  '
$ws.Range('D18').Value = "'92.00"
$ws.Range('E18').Value = '  -2.43%  '
$ws.Range('D19').Value = "'0.06859"
$ws.Range('E19').Value = '  -2.20%  '
$ws.Range('D20').Value = "'18.26"
$ws.Range('E20').Value = '  -7.21%  '
$ws.Range('D21').Value = "'6.595"
$ws.Range('E21').Value = '  -5.73%  '
$ws.Range('B22').Value = 'BitDAO'
$ws.Range('C22').Value = 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'
$ws.Range('D22').Value = "'0.5556"
$ws.Range('E22').Value = '  -6.29%  '
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').Value = "'1.001"
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').Value = "'13.05"
$ws.Range('E24').Value = '  -5.18%  '
$ws.Range('B25').Value = 'WrappedBTC'
$ws.Range('C25').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D25').Value = "'23.110.40"
$ws.Range('E25').Value = '  -3.21%  '
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').Value = "'2.352"
$ws.Range('E26').Value = '  -3.37%  '
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').Value = "'2.762"
$ws.Range('E27').Value = '  -6.88%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = "'21.14"
$ws.Range('E28').Value = '  -4.45%  '
$ws.Range('B29').Value = 'Monero'
$ws.Range('C29').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D29').Value = "'149.87"
$ws.Range('E29').Value = '  -2.77%  '
$ws.Range('B30').Value = 'HuobiToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D30').Value = "'5.263"
$ws.Range('E30').Value = '  -3.17%  '
$ws.Range('B31').Value = 'BitcoinCash'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D31').Value = "'132.77"
$ws.Range('E31').Value = '  -3.96%  '
$ws.Range('B32').Value = 'WEMIXTOKEN'
$ws.Range('C32').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D32').Value = "'2.359"
$ws.Range('E32').Value = '  -5.42%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = "'6.874"
$ws.Range('E33').Value = '  -12.81%  '
$ws.Range('B34').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C34').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D34').Value = "'1.778.53"
$ws.Range('E34').Value = '  -3.09%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').Value = "'0.9584"
$ws.Range('E35').Value = '  -5.70%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').Value = "'0.07710"
$ws.Range('E36').Value = '  -6.02%  '
$ws.Range('B37').Value = 'InternetComputer(DFINITY)'
$ws.Range('C37').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D37').Value = "'6.285"
$ws.Range('E37').Value = '  -5.70%  '
$ws.Range('B38').Value = 'Algorand'
$ws.Range('C38').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D38').Value = "'0.2554"
$ws.Range('E38').Value = '  -4.77%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = "'0.02719"
$ws.Range('E39').Value = '  -6.84%  '
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D40').Value = "'0.08901"
$ws.Range('E40').Value = '  -2.90%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = "'10.12"
$ws.Range('E41').Value = '  -7.03%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = "'1.369"
$ws.Range('E42').Value = '  -3.75%  '
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').Value = "'0.7093"
$ws.Range('E43').Value = '  -6.52%  '
$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').Value = "'12.63"
$ws.Range('E44').Value = '  -7.07%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = "'15.33"
$ws.Range('E45').Value = '  -7.75%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = "'0.6615"
$ws.Range('E46').Value = '  -4.91%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').Value = "'2.318"
$ws.Range('E47').Value = '  -5.65%  '
$ws.Range('B48').Value = 'Frax'
$ws.Range('C48').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D48').Value = "'0.9997"
$ws.Range('E48').Value = '  +0.00%  '
$ws.Range('B49').Value = 'PancakeSwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D49').Value = "'3.996"
$ws.Range('E49').Value = '  -2.63%  '
$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').Value = "'132.11"
$ws.Range('E50').Value = '  -1.62%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = "'0.07940"
$ws.Range('E51').Value = '  -4.43%  '
